$d = $word.ActiveDocument

# 1. Remove the existing _GoBack bookmark (it currently sits right after
#    "while making the cases true." and needs to move to the end of the
#    rewritten sub-goals bullet below).
$oldBm = $d.Bookmarks.Item("_GoBack")
$oldBm.Delete()

# 2. Rewrite the two sub-bullet texts in place. This phrasing also shows
#    up (unchanged) under Problem 3 further down, so only replace the
#    first (Problem 2) occurrence -- wdReplaceOne (1), not wdReplaceAll.
$d.Content.Find.Execute("What are the constraints?", $true, $false, $false, $false, $false, $true, 1, $false, "The constraints are that we need to select the socks in the dark at random and that we need to select the least amount of socks possible.", 1)
$d.Content.Find.Execute("What are the sub-goals?", $true, $false, $false, $false, $false, $true, 1, $false, "The sub-goals are to select one pair of white, one pair of brown, one pair of black, and another pair of any color.", 1)

# 3. Re-create the _GoBack bookmark, collapsed, right after the new
#    sub-goals text (i.e. at the end of that paragraph's content, before
#    the paragraph mark). A collapsed Range sitting exactly at a
#    paragraph's content end confuses Bookmarks.Add directly, so we
#    temporarily insert a placeholder character to push the paragraph
#    boundary out, add the bookmark next to it, then delete the
#    placeholder again; the bookmark stays correctly anchored.
$target = $d.Content
$target.Find.Execute("The sub-goals are to select one pair of white, one pair of brown, one pair of black, and another pair of any color.")
$pos = $target.End

$placeholder = $d.Range($pos, $pos)
$placeholder.InsertAfter("X")

$bmRange = $d.Range($pos, $pos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$placeholderRange = $d.Range($pos, $pos + 1)
$placeholderRange.Delete()
